# LWRS_results_laptop.xlsx — "Update median and median with h2 ptc, fixed
# error in htse cashflows"
#
# The "Median" (row 3) and "Median H2 PTC" (row 4) cases were re-run with
# corrected HTSE cashflow inputs, changing their Mean NPV (B) and Std NPV
# (C) inputs. Every other changed cell on the sheet (D3,E3,F3,G3 and the
# B13:E15 tornado-chart helper block) is a formula that recomputes from
# these four inputs, so we only need to push the four new source values in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated Median (row 3) and Median H2 PTC (row 4) inputs -------------
$ws.Range("B3").Value = 4548180180.3400002
$ws.Range("C3").Value = 6798603.1220399998
$ws.Range("B4").Value = 6178525094.7600002
$ws.Range("C4").Value = 4320401.8213999998

# --- Drop the stale hidden "_xlchart.v1.*" helper names -------------------
# These were auto-generated chart helper defined names left over from an
# earlier chart rebuild; they no longer serve a purpose and are removed.
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "_xlchart.v1.*") {
        $n.Delete()
    }
}

# --- Chart 4 ("Delta NPV") value axis: show in millions, not billions -----
# (Now that the Median is much closer to the baseline, billions is too
# coarse a display unit for the bar.)
$co = $ws.ChartObjects().Item(2)
$valueAxis = $co.Chart.Axes(2)
$valueAxis.DisplayUnit = 6
$valueAxis.HasDisplayUnitLabel = $true

# --- Move the active selection (cosmetic) ----------------------------------
$ws.Range("F20").Select()
